$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "F2" = 1.86
    "G2" = 1.87
    "H2" = 6.6
    "I2" = 6.8
    "J2" = 3.1
    "K2" = 3.2
    "L2" = 0
    "M2" = 0
    "N2" = 0
    "O2" = 0
    "P2" = 5.3
    "Q2" = 1.22
    "R2" = 1.97
    "S2" = 2.02
    "T2" = 0
    "U2" = 0
    "V2" = 1.16
    "W2" = 2.12
    "X2" = 1000
    "Y2" = 1000
    "Z2" = 1000
    "AA2" = 1000
    "AB2" = 1000
    "AC2" = 5.6
    "AD2" = 12
    "AE2" = 55
    "AF2" = 1000
    "AG2" = 4.3
    "AH2" = 9.199999999999999
    "AI2" = 40
    "AJ2" = 1000
    "AK2" = 7.8
    "AL2" = 15.5
    "AM2" = 75
    "AN2" = 11.5
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
